# Updated cryptos list on Fri Oct 13 07:47:41 UTC 2023 with GitHub Actions
# Refresh the "Price" (D) and "Volume(1h)" (E) columns for each coin row
# with the latest scraped values. Price strings can look like plain
# numbers ("61.68"), so the cell is forced to Text format before writing
# so Excel keeps the literal digit-and-dot string instead of re-parsing
# it into a number (which would drop formatting like trailing zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  Price = "26.988.26";  Volume = "  +0.47%  " },
    @{ Row = 3;  Price = "1.555.64";   Volume = "  -0.47%  " },
    @{ Row = 4;  Price = "";           Volume = "  +0.25%  " },
    @{ Row = 5;  Price = "206.88";     Volume = "  +0.48%  " },
    @{ Row = 6;  Price = "";           Volume = "  -0.12%  " },
    @{ Row = 7;  Price = "";           Volume = "  +0.25%  " },
    @{ Row = 8;  Price = "";           Volume = "  +0.41%  " },
    @{ Row = 9;  Price = "21.53";      Volume = "  -1.08%  " },
    @{ Row = 10; Price = "";           Volume = "  -0.45%  " },
    @{ Row = 11; Price = "";           Volume = "  -0.68%  " },
    @{ Row = 12; Price = "1.776.43";   Volume = "  -0.49%  " },
    @{ Row = 13; Price = "1.553.55";   Volume = "  -0.80%  " },
    @{ Row = 14; Price = "";           Volume = "  -0.49%  " },
    @{ Row = 15; Price = "";           Volume = "  -0.35%  " },
    @{ Row = 16; Price = "26.979.05";  Volume = "  +0.39%  " },
    @{ Row = 17; Price = "61.68";      Volume = "  +0.60%  " },
    @{ Row = 18; Price = "214.76";     Volume = "  -0.34%  " },
    @{ Row = 19; Price = "";           Volume = "  +0.62%  " },
    @{ Row = 20; Price = "7.24";       Volume = "  -1.91%  " },
    @{ Row = 21; Price = "";           Volume = "  +0.18%  " },
    @{ Row = 22; Price = "4.05";       Volume = "  -2.19%  " },
    @{ Row = 23; Price = "";           Volume = "  +0.03%  " },
    @{ Row = 24; Price = "";           Volume = "  -2.46%  " },
    @{ Row = 25; Price = "153.85";     Volume = "  -0.14%  " },
    @{ Row = 26; Price = "";           Volume = "  -0.36%  " },
    @{ Row = 27; Price = "14.87";      Volume = "  -0.64%  " },
    @{ Row = 28; Price = "";           Volume = "  +0.24%  " },
    @{ Row = 29; Price = "";           Volume = "  +0.25%  " },
    @{ Row = 30; Price = "";           Volume = "  -0.85%  " },
    @{ Row = 31; Price = "1.10";       Volume = "  -0.67%  " },
    @{ Row = 32; Price = "";           Volume = "  +1.93%  " },
    @{ Row = 33; Price = "1.374.60";   Volume = "  -1.52%  " },
    @{ Row = 34; Price = "2.95";       Volume = "  +1.20%  " },
    @{ Row = 35; Price = "";           Volume = "  +1.46%  " },
    @{ Row = 36; Price = "0.971";      Volume = "  +5.69%  " },
    @{ Row = 37; Price = "";           Volume = "  +0.23%  " },
    @{ Row = 38; Price = "";           Volume = "  +0.36%  " },
    @{ Row = 39; Price = "0.520";      Volume = "  -1.93%  " },
    @{ Row = 40; Price = "0.810";      Volume = "  -0.35%  " },
    @{ Row = 41; Price = "";           Volume = "  +0.33%  " },
    @{ Row = 42; Price = "0.981";      Volume = "  -0.93%  " },
    @{ Row = 43; Price = "";           Volume = "  -0.44%  " },
    @{ Row = 44; Price = "";           Volume = "  +2.14%  " },
    @{ Row = 45; Price = "63.94";      Volume = "  +0.30%  " },
    @{ Row = 46; Price = "1.74";       Volume = "  -2.33%  " },
    @{ Row = 47; Price = "1.689.90";   Volume = "  -0.63%  " },
    @{ Row = 48; Price = "";           Volume = "  -2.96%  " },
    @{ Row = 49; Price = "86.27";      Volume = "  -0.50%  " },
    @{ Row = 50; Price = "";           Volume = "  +0.68%  " },
    @{ Row = 51; Price = "0.0956";     Volume = "  +0.34%  " }
)

foreach ($u in $updates) {
    if ($u.Price -ne "") {
        $priceCell = $ws.Range("D$($u.Row)")
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $u.Price
        $priceCell.Style = "Normal"
    }
    $ws.Range("E$($u.Row)").Value = $u.Volume
}
